# Apply the "Saldo" export refresh: update the rows whose data changed
# (new transactions inserted, some removed, and the report re-sorted by
# descending Saldo). Only the 36 data rows whose contents actually differ
# from the previous export are touched; everything else is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @(3, "004550750", "THEO", 49557.68),
  @(4, "005064129", "THIAGO", 24910.96),
  @(5, "004364200", "BLOCO", 19918.99),
  @(6, "004224405", "MILA", 17710.87),
  @(7, "004690692", "PHYLIA", 14722.73),
  @(8, "004224011", "THOMAS", 14567.86),
  @(9, "005274028", "RAFAEL", 9597.5),
  @(10, "004467884", "ANA", 2612.52),
  @(11, "004369172", "LUIZA", 1161.12),
  @(12, "004487140", "VALMIR", 1160),
  @(13, "004748761", "MARCELO", 1000),
  @(14, "004474776", "GILSON", 973.01),
  @(15, "004392159", "RODRIGO", 900.21),
  @(16, "004855570", "LUISA", 866.25),
  @(17, "005046919", "MARIANA", 800),
  @(19, "001761119", "BLUEMETRIX", 433.11),
  @(20, "004547722", "MARCIA", 400),
  @(42, "004212132", "JOAO", 86.38),
  @(43, "004207278", "CESAR", 84.93),
  @(44, "004239387", "LUIZ", 83.93),
  @(45, "004261201", "ANA", 83.09),
  @(46, "003115072", "VICTOR", 81.27),
  @(47, "005348011", "TATIANA", 80.91),
  @(48, "005009947", "VERANICE", 80.79),
  @(49, "004318604", "RENAN", 80.51),
  @(50, "004994036", "BALTASAR", 80.5),
  @(51, "000330949", "RENATO", 80.09),
  @(52, "004809902", "PEDRO", 79.88),
  @(53, "004267976", "E3", 79.84),
  @(54, "005032151", "ANA", 79.27),
  @(55, "004454365", "RAFAEL", 79.25),
  @(56, "005256849", "SANDRO", 77.17),
  @(57, "004479734", "RODRIGO", 76),
  @(58, "004453045", "JULIAN", 75.83),
  @(59, "004230529", "LAIS", 75.09),
  @(60, "004432579", "ANA", 73.71)
)

foreach ($u in $updates) {
  $r = $u[0]

  # Columns A/B hold account numbers and names and must stay text (the
  # account numbers have significant leading zeros) even though they look
  # numeric. Stamp a text format just for the write so Excel doesn't
  # auto-convert "004550750" -> 4550750, then drop back to General so the
  # cell's number format matches every other untouched row in the sheet.
  $ws.Cells.Item($r, 1).NumberFormat = "@"
  $ws.Cells.Item($r, 1).Value = $u[1]
  $ws.Cells.Item($r, 1).NumberFormat = "General"

  $ws.Cells.Item($r, 2).NumberFormat = "@"
  $ws.Cells.Item($r, 2).Value = $u[2]
  $ws.Cells.Item($r, 2).NumberFormat = "General"

  $ws.Cells.Item($r, 3).Value = $u[3]
}
